$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows for "RM 232" (row 26) and "SC 92" (row 28, which becomes row 27
# after the first deletion shifts rows up). Deleting shifts remaining rows up,
# matching the new dimension A1:F33.
$ws.Rows("26:26").Delete()
$ws.Rows("27:27").Delete()

# Apply the remaining per-cell value updates (imputed/removed values) to match
# the target missing-data pattern.
$ws.Range("C2").Value = 14.9
$ws.Range("F2").Value = 18.03
$ws.Range("D3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("D4").Value = -15.4
$ws.Range("D5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("C12").Value = 12.5
$ws.Range("F13").Value = 17.1
$ws.Range("C14").ClearContents()
$ws.Range("F19").ClearContents()
$ws.Range("C20").Value = 12.5
$ws.Range("C21").Value = 12.7
$ws.Range("C22").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("D23").Value = -13.9
$ws.Range("F25").Value = 16.6
$ws.Range("D27").ClearContents()
$ws.Range("F28").Value = 17.44
$ws.Range("D29").Value = -13
$ws.Range("B30").Value = -19.7
$ws.Range("C31").Value = 15.3
$ws.Range("F31").ClearContents()
$ws.Range("B32").ClearContents()
$ws.Range("F32").Value = 17.39
$ws.Range("C33").Value = 10.4
